# edit.ps1 - applies the "Engagement Scope" / "Investment Summary" table
# restructuring + title-slide date bump described by the commit diff.

function Hex2Bgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$COLOR_HEADER_FILL = Hex2Bgr "A01C02"
$COLOR_HEADER_TEXT = Hex2Bgr "FFFFFF"
$COLOR_DATA_FILL   = Hex2Bgr "E7E6E6"
$COLOR_SPACER_FILL = Hex2Bgr "FFFFFF"
$COLOR_SPACER_LINE = Hex2Bgr "D0D0D0"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Title slide: bump the date in the byline.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item(5).TextFrame.TextRange.Text = "Alison Smith | November 17, 2025"

# ---------------------------------------------------------------------------
# 2. Slide 3 "Engagement Scope" table: collapse the 3-col / 17-row table
#    (Category | Parameter | Scope) into a 5-col / 9-row table
#    (Parameter | Scope | spacer | Parameter | Scope), pairing data row N
#    with data row N+8, bolding the parameter labels.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$tbl1 = $slide3.Shapes.Item(3).Table

# Right-hand pairing data (was rows 10..17 / old data rows 9..16).
$rightData = @(
  @("Deployment Regions", "3 Azure regions globally"),
  @("Availability Requirements", "High availability (99.95%)"),
  @("Infrastructure Complexity", "vWAN + Azure Firewall + routing"),
  @("Security Requirements", "Azure Firewall DDoS protection"),
  @("Compliance Frameworks", "SOC2 ISO27001"),
  @("Latency Requirements", "<50ms inter-region"),
  @("Routing Complexity", "BGP route propagation"),
  @("Deployment Environments", "2 environments (non-prod prod)")
)

# Drop the "Category" column (col 1); Parameter/Scope remain as cols 1-2.
$tbl1.Columns.Item(1).Delete()

# Drop the old right-hand data rows (rows 10-17 after the column delete);
# their text has been captured above in $rightData for the new layout.
for ($i = 17; $i -ge 10; $i--) {
    $tbl1.Rows.Item($i).Delete()
}

# Add spacer + Parameter + Scope columns (cols 3, 4, 5).
$tbl1.Columns.Add() | Out-Null
$tbl1.Columns.Add() | Out-Null
$tbl1.Columns.Add() | Out-Null

# Column widths (EMU / 12700 = points).
$tbl1.Columns.Item(1).Width = 1567967 / 12700.0
$tbl1.Columns.Item(2).Width = 2526170 / 12700.0
$tbl1.Columns.Item(3).Width = 435546 / 12700.0
$tbl1.Columns.Item(4).Width = 1567967 / 12700.0
$tbl1.Columns.Item(5).Width = 2613279 / 12700.0

# --- Header row (row 1) ---
$headerTexts = @("Parameter", "Scope", "", "Parameter", "Scope")
for ($c = 1; $c -le 5; $c++) {
    $cell = $tbl1.Cell(1, $c)
    if ($c -eq 3) {
        $cell.Shape.Fill.ForeColor.RGB = $COLOR_SPACER_FILL
        $lnL = $cell.Borders(2)
        $lnL.ForeColor.RGB = $COLOR_SPACER_LINE
        $lnL.Weight = 1
        $lnL.Visible = -1
    } else {
        $cell.Shape.Fill.ForeColor.RGB = $COLOR_HEADER_FILL
        $tr = $cell.Shape.TextFrame.TextRange
        $tr.Text = $headerTexts[$c - 1]
        $tr.Font.Size = 14
        $tr.Font.Color.RGB = $COLOR_HEADER_TEXT
    }
}

# --- Data rows (rows 2-9) ---
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $right = $rightData[$i]

    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl1.Cell($row, $c)
        if ($c -eq 3) {
            $cell.Shape.Fill.ForeColor.RGB = $COLOR_SPACER_FILL
            $lnL = $cell.Borders(2)
            $lnL.ForeColor.RGB = $COLOR_SPACER_LINE
            $lnL.Weight = 1
            $lnL.Visible = -1
            continue
        }

        $cell.Shape.Fill.ForeColor.RGB = $COLOR_DATA_FILL
        $tr = $cell.Shape.TextFrame.TextRange

        if ($c -eq 4) { $tr.Text = $right[0] }
        elseif ($c -eq 5) { $tr.Text = $right[1] }

        $tr.Font.Size = 11
        if ($c -eq 1 -or $c -eq 4) {
            $tr.Font.Bold = 1
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Slide 9 "Investment Summary" table: drop the 4 zero-value line-item
#    rows, keeping the header row and the TOTAL INVESTMENT row.
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$tbl2 = $slide9.Shapes.Item(3).Table

for ($i = 5; $i -ge 2; $i--) {
    $tbl2.Rows.Item($i).Delete()
}
